$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set the value for F2 (previously empty)
$ws.Range("F2").Value = 2

# Force recalculation so H2's SUM formula reflects the new value
$excel.Calculate()

# Update the active selection to F3 to match the saved cursor position
$ws.Range("F3").Select()
